$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 4-11 (the "groups") get the Country value ("United States") that the
# top two rows already have, instead of "unknown".
$ws.Range("F4").Value = "United States"
$ws.Range("F5").Value = "United States"
$ws.Range("F6").Value = "United States"
$ws.Range("F7").Value = "United States"
$ws.Range("F8").Value = "United States"
$ws.Range("F9").Value = "United States"
$ws.Range("F10").Value = "United States"
$ws.Range("F11").Value = "United States"

# Turn duplicate IP addresses (...26) into the already-used "groups" (...27),
# along with their matching Reverse DNS hostnames.
$ws.Range("B8").Value = "142.250.123.27"
$ws.Range("D8").Value = "gh-in-f27.1e100.net"

$ws.Range("B10").Value = "108.177.12.27"
$ws.Range("D10").Value = "ua-in-f27.1e100.net"

$ws.Range("B11").Value = "64.233.186.27"
$ws.Range("D11").Value = "cb-in-f27.1e100.net"
